$wb = $excel.ActiveWorkbook

# Update employment-share targets (fracEmployed) across the employment_* sheets.
# Values below are the new simulated-employment-share targets replacing the
# previous isEmployed-boolean-derived targets, per "Apply employment share
# per BU rather than boolean isEmployed in employment alignment".

# --- employment_acfemales ---
$ws = $wb.Worksheets.Item("employment_acfemales")
$ws.Range("B2").Value = 0.38052187841351931
$ws.Range("B3").Value = 0.39494946229994121
$ws.Range("B4").Value = 0.3708155540743337
$ws.Range("B5").Value = 0.39400614231074099
$ws.Range("B6").Value = 0.39038070193201418
$ws.Range("B7").Value = 0.4045626857325601
$ws.Range("B8").Value = 0.43452589116520407
$ws.Range("B9").Value = 0.41769307045724119
$ws.Range("B10").Value = 0.42367705016359092
$ws.Range("B11").Value = 0.45367977278524363
$ws.Range("B12").Value = 0.44799750636445951
$ws.Range("B13").Value = 0.43338321233616106
$ws.Range("B14").Value = 0.43184882342437486

# --- employment_acmales ---
$ws = $wb.Worksheets.Item("employment_acmales")
$ws.Range("B2").Value = 0.51193623756613316
$ws.Range("B3").Value = 0.51388796407993187
$ws.Range("B4").Value = 0.50731054213776516
$ws.Range("B5").Value = 0.51131531460632029
$ws.Range("B6").Value = 0.53318545452592403
$ws.Range("B7").Value = 0.58285048457054911
$ws.Range("B8").Value = 0.59181134292156978
$ws.Range("B9").Value = 0.56828787809135639
$ws.Range("B10").Value = 0.57210350926861453
$ws.Range("B11").Value = 0.56859266378921847
$ws.Range("B12").Value = 0.61407456424300411
$ws.Range("B13").Value = 0.60188645030954646
$ws.Range("B14").Value = 0.59164974063354814

# --- employment_femalewdep ---
$ws = $wb.Worksheets.Item("employment_femalewdep")
$ws.Range("B2").Value = 0.30785190812656721
$ws.Range("B3").Value = 0.30811218783286087
$ws.Range("B4").Value = 0.31195865499437081
$ws.Range("B5").Value = 0.28581726919907879
$ws.Range("B6").Value = 0.27963714773111775
$ws.Range("B7").Value = 0.30752858106411018
$ws.Range("B8").Value = 0.30346469681990601
$ws.Range("B9").Value = 0.31364752176359245
$ws.Range("B10").Value = 0.35381764642659308
$ws.Range("B11").Value = 0.34820741422940238
$ws.Range("B12").Value = 0.3607173104243539
$ws.Range("B13").Value = 0.36274965739361886
$ws.Range("B14").Value = 0.36845418733118185

# --- employment_malewdep ---
$ws = $wb.Worksheets.Item("employment_malewdep")
$ws.Range("B2").Value = 0.41894661081369583
$ws.Range("B3").Value = 0.40321361657827082
$ws.Range("B4").Value = 0.39961607117918457
$ws.Range("B5").Value = 0.39812522211769624
$ws.Range("B6").Value = 0.40991729842040864
$ws.Range("B7").Value = 0.40545695127308695
$ws.Range("B8").Value = 0.41228362735812552
$ws.Range("B9").Value = 0.42686713226978734
$ws.Range("B10").Value = 0.42517186848982991
$ws.Range("B11").Value = 0.41013480515252171
$ws.Range("B12").Value = 0.43258856779953725
$ws.Range("B13").Value = 0.44326586886465918
$ws.Range("B14").Value = 0.46207239052048033

# --- employment_smales ---
$ws = $wb.Worksheets.Item("employment_smales")
$ws.Range("B2").Value = 0.44580390908196338
$ws.Range("B3").Value = 0.44324938187318652
$ws.Range("B4").Value = 0.44437298985301549
$ws.Range("B5").Value = 0.45665144279125908
$ws.Range("B6").Value = 0.47060122836212664
$ws.Range("B7").Value = 0.46283571902209464
$ws.Range("B8").Value = 0.4829455027595948
$ws.Range("B9").Value = 0.49083580436737168
$ws.Range("B10").Value = 0.48680351073667416
$ws.Range("B11").Value = 0.49651737484706926
$ws.Range("B12").Value = 0.52320974701450895
$ws.Range("B13").Value = 0.53678039526990218
$ws.Range("B14").Value = 0.52108578146447826

# --- employment_sfemales ---
$ws = $wb.Worksheets.Item("employment_sfemales")
$ws.Range("B2").Value = 0.30208181732795336
$ws.Range("B3").Value = 0.30061108532233771
$ws.Range("B4").Value = 0.30276407721188736
$ws.Range("B5").Value = 0.31145040929785739
$ws.Range("B6").Value = 0.31244190895581908
$ws.Range("B7").Value = 0.31288817157553944
$ws.Range("B8").Value = 0.30967743686400945
$ws.Range("B9").Value = 0.3052372997493929
$ws.Range("B10").Value = 0.31519769527800484
$ws.Range("B11").Value = 0.31410432759222856
$ws.Range("B12").Value = 0.319650585035544
$ws.Range("B13").Value = 0.34466439690323852
$ws.Range("B14").Value = 0.34069381864858117

# --- employment_couples ---
$ws = $wb.Worksheets.Item("employment_couples")
$ws.Range("B2").Value = 0.82383029192035362
$ws.Range("B3").Value = 0.821731322462114
$ws.Range("B4").Value = 0.81609950788251906
$ws.Range("B5").Value = 0.83106676793864132
$ws.Range("B6").Value = 0.83335242550016386
$ws.Range("B7").Value = 0.83973543797809402
$ws.Range("B8").Value = 0.84157973500650807
$ws.Range("B9").Value = 0.84545298666168944
$ws.Range("B10").Value = 0.8534134938695962
$ws.Range("B11").Value = 0.86500591632281076
$ws.Range("B12").Value = 0.88091732038400394
$ws.Range("B13").Value = 0.87862197711765877
$ws.Range("B14").Value = 0.89231035103286582

# The active sheet moves from "employment_smales" to "employment_acmales"
# (bookViews firstSheet/activeTab + per-sheet tabSelected in the XML).
$wsActive = $wb.Worksheets.Item("employment_acmales")
$wsActive.Activate()
